$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 becomes a real number instead of the inline text "0987654321"
$ws.Range("B2").Value = 987654321

# New rows 3 and 4 (the removed "hospital" class rows)
$ws.Range("A3").Value = "apollo@prac.to"
$ws.Range("C3").Value = "Bangalore"

$ws.Range("A4").Value = "practo@oafaz.com"
$ws.Range("C4").Value = "snfkjad"

# B3 / B4 stay as text (not auto-converted numbers) like the other string
# cells in the sheet, so force a text format before assigning, then drop
# back to the default style afterwards.
$ws.Range("B3:B4").NumberFormat = "@"
$ws.Range("B3").Value = "1234567890"
$ws.Range("B4").Value = "1234567890"
$ws.Range("B3:B4").Style = "Normal"
